$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Locate the "PWS Number" row within the first table (column 1 contains the
# field name). This is more robust than hard-coding a row index.
$targetRow = -1
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $label = $tbl.Cell($i, 1).Range.Text
    if ($label -match "^PWS Number") {
        $desc = $tbl.Cell($i, 2).Range.Text
        if ($desc -notmatch "MUST be the same value") {
            $targetRow = $i
            break
        }
    }
}

if ($targetRow -eq -1) {
    throw "Could not locate target PWS Number row"
}

# --- Edit 1: Input Message cell - split the description into two
#     paragraphs and reword the parenthetical requirement text. ---
$cell2 = $tbl.Cell($targetRow, 2)
$p2 = $cell2.Range.Paragraphs.Item(1)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>
<w:p w14:paraId="006D8F6F" w14:textId="4A2EFB61" w:rsidR="00CE5C05" w:rsidRPr="00CE5C05" w:rsidRDefault="00CE5C05" w:rsidP="00CE5C05"><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00CE5C05"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>State-assigned Public Water System identifier.</w:t></w:r></w:p>
<w:p><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>(Required, must be 9 characters.)</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r2.InsertXML($xml2)

# --- Edit 2: Error Message cell - reword "Must not be longer than 9"
#     into "Must be exactly 9" (split across runs per the source edit). ---
$cell3 = $tbl.Cell($targetRow, 3)
$p3 = $cell3.Range.Paragraphs.Item(1)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>
<w:p w14:paraId="06243C37" w14:textId="0D18B235" w:rsidR="00CE5C05" w:rsidRPr="00CE5C05" w:rsidRDefault="00CE5C05" w:rsidP="00CE5C05"><w:pPr><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="00CE5C05"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">Must </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>be exactly</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="005A68F4"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">9 </w:t></w:r><w:r w:rsidRPr="00CE5C05"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>characters.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r3.InsertXML($xml3)

Write-Host "Cell2 text:" $tbl.Cell($targetRow, 2).Range.Text
Write-Host "Cell3 text:" $tbl.Cell($targetRow, 3).Range.Text
